$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: nurse/doctor (unchanged) -> enter patient's vitals
$ws.Range("C5").Value = "enter patient's vitals"

# Row 6: doctor -> enter patient's vitals, prescriptions, and lab tests
$ws.Range("B6").Value = "doctor"
$ws.Range("C6").Value = "enter patient's vitals, prescriptions, and lab tests"

# Row 7: patient -> enter personal information
$ws.Range("B7").Value = "patient"
$ws.Range("C7").Value = "enter personal information"

# "so that..." column (D) updated for rows 5-7 to the new shared note
$ws.Range("D5").Value = "the system can be udpated to view the information"
$ws.Range("D6").Value = "the system can be udpated to view the information"
$ws.Range("D7").Value = "the system can be udpated to view the information"

# Row heights for rows 5-7 changed to 45
$ws.Rows("5").RowHeight = 45
$ws.Rows("6").RowHeight = 45
$ws.Rows("7").RowHeight = 45

# Selection changed to E6
$ws.Range("E6").Select()
